$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 3032.7058
$ws.Range("J116").Value = 4169.5
$ws.Range("L116").Value = 4169.5
$ws.Range("N116").Value = -11053.5

$ws.Range("H138").Value = 4362.9194
$ws.Range("I138").Value = 6094.654
$ws.Range("J138").Value = 3112.2222
$ws.Range("K138").Value = 18283.962
$ws.Range("L138").Value = 9336.6666
$ws.Range("M138").Value = -13143.962
$ws.Range("N138").Value = -19616.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 8065.2
$ws.Range("I37").Value = 4666.3335
$ws.Range("J37").Value = 9521.857
$ws.Range("K37").Value = 4666.3335
$ws.Range("L37").Value = 9521.857
$ws.Range("M37").Value = -4393.3335
$ws.Range("N37").Value = -10067.857

$ws.Range("H44").Value = 21574.5
$ws.Range("J44").Value = 21574.5
$ws.Range("L44").Value = 21574.5
$ws.Range("N44").Value = -22550.5

$ws.Range("H49").Value = 30000
$ws.Range("J49").Value = 30000
$ws.Range("L49").Value = 30000
$ws.Range("N49").Value = -30520

$ws.Range("H55").Value = 17775.666
$ws.Range("J55").Value = 17775.666
$ws.Range("L55").Value = 17775.666
$ws.Range("N55").Value = -18405.666

$ws.Range("H63").Value = 4161
$ws.Range("I63").Value = 1601.6666
$ws.Range("J63").Value = 8000
$ws.Range("K63").Value = 1601.6666
$ws.Range("L63").Value = 8000
$ws.Range("M63").Value = -915.6666
$ws.Range("N63").Value = -9372

$ws.Range("H66").Value = 4161
$ws.Range("I66").Value = 1601.6666
$ws.Range("J66").Value = 8000
$ws.Range("K66").Value = 8008.333000000001
$ws.Range("L66").Value = 40000
$ws.Range("M66").Value = -4576.333000000001
$ws.Range("N66").Value = -46864

$ws.Range("H80").Value = 21960
$ws.Range("J80").Value = 21960
$ws.Range("L80").Value = 21960
$ws.Range("N80").Value = -23956

$ws.Range("H83").Value = 21960
$ws.Range("J83").Value = 21960
$ws.Range("L83").Value = 65880
$ws.Range("N83").Value = -75864

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 17550
$ws.Range("J35").Value = 17550
$ws.Range("L35").Value = 17550
$ws.Range("N35").Value = -18170

$ws.Range("H82").Value = 17024.273
$ws.Range("I82").Value = 10125
$ws.Range("J82").Value = 20966.715
$ws.Range("K82").Value = 10125
$ws.Range("L82").Value = 20966.715
$ws.Range("M82").Value = -9742
$ws.Range("N82").Value = -21732.715

$ws.Range("H85").Value = 17024.273
$ws.Range("I85").Value = 10125
$ws.Range("J85").Value = 20966.715
$ws.Range("K85").Value = 10125
$ws.Range("L85").Value = 20966.715
$ws.Range("M85").Value = -8799
$ws.Range("N85").Value = -23618.715

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 14900
$ws.Range("J112").Value = 14900
$ws.Range("L112").Value = 14900
$ws.Range("N112").Value = -17854

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 141411410
$ws.Range("J96").Value = 141411410
$ws.Range("L96").Value = 424234230
$ws.Range("N96").Value = -424238348

$ws.Range("H98").Value = 507.36365
$ws.Range("I98").Value = 431.83334
$ws.Range("K98").Value = 1295.50002
$ws.Range("M98").Value = 202.4999800000001

$ws.Range("H110").Value = 3150
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H131").Value = 3070.413
$ws.Range("I131").Value = 10280
$ws.Range("J131").Value = 2191.195
$ws.Range("K131").Value = 30840
$ws.Range("L131").Value = 6573.585000000001
$ws.Range("M131").Value = -25800
$ws.Range("N131").Value = -16653.585

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 940
$ws.Range("I43").Value = 940
$ws.Range("K43").Value = 940
$ws.Range("M43").Value = -789

$ws.Range("H46").Value = 17950
$ws.Range("J46").Value = 17950
$ws.Range("L46").Value = 17950
$ws.Range("N46").Value = -18262

$ws.Range("H80").Value = 102628.27
$ws.Range("I80").Value = 3425.625
$ws.Range("J80").Value = 367168.66
$ws.Range("K80").Value = 3425.625
$ws.Range("L80").Value = 367168.66
$ws.Range("M80").Value = -2427.625
$ws.Range("N80").Value = -369164.66

$ws.Range("H83").Value = 102628.27
$ws.Range("I83").Value = 3425.625
$ws.Range("J83").Value = 367168.66
$ws.Range("K83").Value = 17128.125
$ws.Range("L83").Value = 1835843.3
$ws.Range("M83").Value = -12136.125
$ws.Range("N83").Value = -1845827.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H94").Value = 25000
$ws.Range("J94").Value = 25000
$ws.Range("L94").Value = 25000
$ws.Range("N94").Value = -26352

$ws.Range("H136").Value = 2083.6316
$ws.Range("I136").Value = 1199.5714
$ws.Range("J136").Value = 4559
$ws.Range("K136").Value = 3598.7142
$ws.Range("L136").Value = 13677
$ws.Range("M136").Value = -1048.7142
$ws.Range("N136").Value = -18777

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 14207.429
$ws.Range("I54").Value = 9750
$ws.Range("J54").Value = 15990.4
$ws.Range("K54").Value = 9750
$ws.Range("L54").Value = 15990.4
$ws.Range("M54").Value = -9230
$ws.Range("N54").Value = -17030.4

$ws.Range("H81").Value = 4515.5454
$ws.Range("I81").Value = 2223.6667
$ws.Range("J81").Value = 5375
$ws.Range("K81").Value = 4447.3334
$ws.Range("L81").Value = 10750
$ws.Range("M81").Value = -3386.3334
$ws.Range("N81").Value = -12872

$ws.Range("H84").Value = 4515.5454
$ws.Range("I84").Value = 2223.6667
$ws.Range("J84").Value = 5375
$ws.Range("K84").Value = 22236.667
$ws.Range("L84").Value = 53750
$ws.Range("M84").Value = -16932.667
$ws.Range("N84").Value = -64358

$ws.Range("H132").Value = 3601.5557
$ws.Range("I132").Value = 4049.158
$ws.Range("J132").Value = 2538.5
$ws.Range("K132").Value = 12147.474
$ws.Range("L132").Value = 7615.5
$ws.Range("M132").Value = -9617.474
$ws.Range("N132").Value = -12675.5

$ws.Range("H136").Value = 7816.483
$ws.Range("I136").Value = 9477.087
$ws.Range("J136").Value = 1450.8334
$ws.Range("K136").Value = 28431.261
$ws.Range("L136").Value = 4352.5002
$ws.Range("M136").Value = -25881.261
$ws.Range("N136").Value = -9452.5002
